$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates that are plain text (kept as text automatically) ---
$dPlainRows = @(2, 3, 7, 13, 16, 18, 24, 33, 50)
$dPlainVals = @("65.856.32", "3.607.29", "3.607.61", "4.221.84", "3.601.99", "65.939.56", "3.753.67", "3.614.07", "2.468.26")
for ($i = 0; $i -lt $dPlainRows.Length; $i++) {
    $ws.Cells.Item($dPlainRows[$i], 4).Value = $dPlainVals[$i]
}

# --- Price (column D) updates that look like plain numbers; force them to stay text ---
$dNumRows = @(4, 5, 6, 12, 14, 15, 19, 20, 21, 22, 28, 29, 37, 41, 42, 43, 44, 45, 48)
$dNumVals = @("1.00", "605.49", "137.49", "0.393", "28.13", "0.0000188", "10.13", "14.72", "5.94", "399.09", "8.24", "1.69", "5.40", "0.0839", "0.843", "26.15", "43.45", "1.26", "1.72")
for ($i = 0; $i -lt $dNumRows.Length; $i++) {
    $cell = $ws.Cells.Item($dNumRows[$i], 4)
    $cell.NumberFormat = "@"
    $cell.Value = $dNumVals[$i]
    $cell.Style = "Normal"
}

# --- Volume(1h) (column E) updates; values always contain "%" so they remain text ---
$eRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 27, 28, 29, 30, 31, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)
$eVals = @("  +0.88%  ", "  +1.87%  ", "  +0.01%  ", "  +0.35%  ", "  -2.00%  ", "  +1.85%  ", "  +0.04%  ", "  +0.78%  ", "  +0.25%  ", "  -0.06%  ", "  +2.03%  ", "  +3.15%  ", "  -0.10%  ", "  +1.34%  ", "  -0.34%  ", "  +0.85%  ", "  -2.22%  ", "  +2.58%  ", "  -0.23%  ", "  +0.67%  ", "  +3.00%  ", "  +1.95%  ", "  +0.62%  ", "  +2.02%  ", "  +4.71%  ", "  +29.84%  ", "  +4.42%  ", "  +4.58%  ", "  +1.86%  ", "  +3.66%  ", "  +0.89%  ", "  -0.01%  ", "  +7.85%  ", "  +2.93%  ", "  +1.33%  ", "  +1.03%  ", "  +2.70%  ", "  +1.43%  ", "  -0.77%  ", "  +1.17%  ", "  +4.02%  ", "  +2.39%  ", "  +0.02%  ", "  +1.35%  ", "  +4.04%  ", "  +0.58%  ", "  +3.36%  ")
for ($i = 0; $i -lt $eRows.Length; $i++) {
    $ws.Cells.Item($eRows[$i], 5).Value = $eVals[$i]
}

